$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Team member list: insert last names
# ---------------------------------------------------------------
$d.Content.Find.Execute("Ji-su", $true, $false, $false, $false, $false, $true, 1, $false, "Ji-su Choi", 2) | Out-Null
$d.Content.Find.Execute("Jinho", $true, $false, $false, $false, $false, $true, 1, $false, "Jinho Bae", 2) | Out-Null
$d.Content.Find.Execute("Tyler", $true, $false, $false, $false, $false, $true, 1, $false, "Tyler McKerihan", 2) | Out-Null
$d.Content.Find.Execute("Aiden", $true, $false, $false, $false, $false, $true, 1, $false, "Aiden Lamb", 2) | Out-Null

# ---------------------------------------------------------------
# 2. Extend the "Facebook Messenger" sentence
# ---------------------------------------------------------------
$d.Content.Find.Execute("Facebook Messenger.", $true, $false, $false, $false, $false, $true, 1, $false, "Facebook Messenger, as well as GitHub READMEs and commit descriptions.", 2) | Out-Null

Write-Host "done phase 1"

# ---------------------------------------------------------------
# 3. Remove the manual line break at the end of the "We aimed to
#    deliver..." paragraph and replace it with a new sentence.
# ---------------------------------------------------------------
$d.Content.Find.Execute("criticism.^l", $true, $false, $false, $false, $false, $true, 1, $false, "criticism. We frequently engaged in group discussions about the direction of the product as a whole as well as individual elements.", 2) | Out-Null

# ---------------------------------------------------------------
# 4. Extend the "We used GitHub..." paragraph with a new sentence.
# ---------------------------------------------------------------
$d.Content.Find.Execute("inexperience using the platform.", $true, $false, $false, $false, $false, $true, 1, $false, "inexperience using the platform. Another related issue that we had was the use of absolute paths in some of our code, leading to errors when used on other systems. We alleviated this issue by using relative paths.", 2) | Out-Null

Write-Host "done phase 2"

# ---------------------------------------------------------------
# 5. Remove the stray "_GoBack" bookmark from the end of the
#    "...demonstrate all of the workshop homework on time." paragraph
#    (it is going to be relocated to the new final paragraph below).
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------
# 6. Insert the new closing paragraph ("As a whole, ...") after the
#    "We aimed to deliver..." paragraph, carrying the relocated
#    "_GoBack" bookmark plus a manual line break and a tab.
# ---------------------------------------------------------------
$anchorRange = $d.Content
$anchorRange.Find.Execute("We frequently engaged in group discussions about the direction of the product as a whole as well as individual elements.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchorRange.Paragraphs(1)
$anchorPara.Range.InsertParagraphAfter()

$newParaIndex = $anchorPara.Index + 1
$newPara = $d.Paragraphs($newParaIndex)
$newRange = $newPara.Range.Duplicate
$newRange.Collapse(1)

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$closingXml = "<w:p $ns>" +
    "<w:r><w:t>As a whole, we are quite satisfied with our performance in Sprint 1 and look forward to employing similar tactics and practises going into Sprint 2.</w:t></w:r>" +
    "<w:bookmarkStart w:id=`"500`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"500`"/>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:tab/></w:r>" +
    "</w:p>"
$newRange.InsertXML($closingXml)

Write-Host "done phase 3"


